# Insert a new data row at row 107 (shifts existing rows 107-196 down to 108-197)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 107, shifting cells down and
# carrying the formatting from the row above (default Excel behaviour).
$ws.Rows(107).Insert()

# New record values (commit: "Fruta / hortaliza, semanal")
$ws.Range("A107").Value = 8
$ws.Range("B107").Value = "Terminal La Palmera de La Serena"
$ws.Range("C107").Value = "Coquimbo"
$ws.Range("D107").Value = [DateTime]"2023-11-03"
$ws.Range("E107").Value = 4
$ws.Range("F107").Value = 100114007
$ws.Range("G107").Value = "Jengibre"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 27000
$ws.Range("L107").Value = 28000
$ws.Range("M107").Value = 27500
$ws.Range("N107").Value = "`$/caja 13 kilos"
$ws.Range("O107").Value = "Perú"
$ws.Range("P107").Value = 2115
$ws.Range("Q107").Value = 13
$ws.Range("R107").Value = "Hortaliza"
